$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 7).Value = 25.35940266666667
$ws.Cells.Item(2, 8).Value = 76.078208
$ws.Cells.Item(2, 9).Value = 0.005186643687654987
$ws.Cells.Item(2, 10).Value = 0.005186643687654986
$ws.Cells.Item(2, 13).Value = 3.425446666666666
$ws.Cells.Item(2, 14).Value = 10.27634
$ws.Cells.Item(2, 15).Value = 0.6657953389778073
$ws.Cells.Item(2, 16).Value = 0.6657953389778073
$ws.Cells.Item(2, 17).Value = 86.8672813331911
$ws.Cells.Item(2, 18).Value = 781.80553199872
$ws.Cells.Item(2, 19).Value = 0.003453243192179357
$ws.Cells.Item(2, 20).Value = 0.003453243192179356
$ws.Cells.Item(3, 7).Value = 25.35940266666667
$ws.Cells.Item(3, 8).Value = 76.078208
$ws.Cells.Item(3, 9).Value = 0.005186643687654987
$ws.Cells.Item(3, 10).Value = 0.005186643687654986
$ws.Cells.Item(3, 15).Value = 0.2094791321596951
$ws.Cells.Item(3, 16).Value = 0.2094791321596952
$ws.Cells.Item(3, 17).Value = 27.33104550519467
$ws.Cells.Item(3, 18).Value = 245.979409546752
$ws.Cells.Item(3, 19).Value = 0.001086493618511528
$ws.Cells.Item(3, 20).Value = 0.001086493618511528
$ws.Cells.Item(4, 7).Value = 25.35940266666667
$ws.Cells.Item(4, 8).Value = 76.078208
$ws.Cells.Item(4, 9).Value = 0.005186643687654987
$ws.Cells.Item(4, 10).Value = 0.005186643687654986
$ws.Cells.Item(4, 13).Value = 0.62317
$ws.Cells.Item(4, 14).Value = 1.86951
$ws.Cells.Item(4, 15).Value = 0.1211239647746572
$ws.Cells.Item(4, 16).Value = 0.1211239647746572
$ws.Cells.Item(4, 17).Value = 15.80321895978667
$ws.Cells.Item(4, 18).Value = 142.22897063808
$ws.Cells.Item(4, 19).Value = 0.0006282268473222207
$ws.Cells.Item(4, 20).Value = 0.0006282268473222206
$ws.Cells.Item(5, 7).Value = 25.35940266666667
$ws.Cells.Item(5, 8).Value = 76.078208
$ws.Cells.Item(5, 9).Value = 0.005186643687654987
$ws.Cells.Item(5, 10).Value = 0.005186643687654986
$ws.Cells.Item(5, 13).Value = 0.01852966666666667
$ws.Cells.Item(5, 14).Value = 0.055589
$ws.Cells.Item(5, 15).Value = 0.003601564087840353
$ws.Cells.Item(5, 16).Value = 0.003601564087840353
$ws.Cells.Item(5, 17).Value = 0.4699012782791112
$ws.Cells.Item(5, 18).Value = 4.229111504512
$ws.Cells.Item(5, 19).Value = 0.00001868002964188206
$ws.Cells.Item(5, 20).Value = 0.00001868002964188206
$ws.Cells.Item(6, 9).Value = 0.9837462940761621
$ws.Cells.Item(6, 10).Value = 0.983746294076162
$ws.Cells.Item(6, 13).Value = 3.425446666666666
$ws.Cells.Item(6, 14).Value = 10.27634
$ws.Cells.Item(6, 15).Value = 0.6657953389778073
$ws.Cells.Item(6, 16).Value = 0.6657953389778073
$ws.Cells.Item(6, 17).Value = 16476.04332092353
$ws.Cells.Item(6, 18).Value = 148284.3898883118
$ws.Cells.Item(6, 19).Value = 0.6549736973326
$ws.Cells.Item(6, 20).Value = 0.6549736973326
$ws.Cells.Item(7, 9).Value = 0.9837462940761621
$ws.Cells.Item(7, 10).Value = 0.983746294076162
$ws.Cells.Item(7, 15).Value = 0.2094791321596951
$ws.Cells.Item(7, 16).Value = 0.2094791321596952
$ws.Cells.Item(7, 19).Value = 0.2060743199483907
$ws.Cells.Item(7, 20).Value = 0.2060743199483907
$ws.Cells.Item(8, 9).Value = 0.9837462940761621
$ws.Cells.Item(8, 10).Value = 0.983746294076162
$ws.Cells.Item(8, 13).Value = 0.62317
$ws.Cells.Item(8, 14).Value = 1.86951
$ws.Cells.Item(8, 15).Value = 0.1211239647746572
$ws.Cells.Item(8, 16).Value = 0.1211239647746572
$ws.Cells.Item(8, 17).Value = 2997.383090565294
$ws.Cells.Item(8, 18).Value = 26976.44781508764
$ws.Cells.Item(8, 19).Value = 0.1191552514708806
$ws.Cells.Item(8, 20).Value = 0.1191552514708806
$ws.Cells.Item(9, 9).Value = 0.9837462940761621
$ws.Cells.Item(9, 10).Value = 0.983746294076162
$ws.Cells.Item(9, 13).Value = 0.01852966666666667
$ws.Cells.Item(9, 14).Value = 0.055589
$ws.Cells.Item(9, 15).Value = 0.003601564087840353
$ws.Cells.Item(9, 16).Value = 0.003601564087840353
$ws.Cells.Item(9, 17).Value = 89.1257755355329
$ws.Cells.Item(9, 18).Value = 802.131979819796
$ws.Cells.Item(9, 19).Value = 0.00354302532429074
$ws.Cells.Item(9, 20).Value = 0.00354302532429074
$ws.Cells.Item(10, 7).Value = 51.27300266666666
$ws.Cells.Item(10, 8).Value = 153.819008
$ws.Cells.Item(10, 9).Value = 0.01048663484403512
$ws.Cells.Item(10, 10).Value = 0.01048663484403512
$ws.Cells.Item(10, 13).Value = 3.425446666666666
$ws.Cells.Item(10, 14).Value = 10.27634
$ws.Cells.Item(10, 15).Value = 0.6657953389778073
$ws.Cells.Item(10, 16).Value = 0.6657953389778073
$ws.Cells.Item(10, 17).Value = 175.6329360745244
$ws.Cells.Item(10, 18).Value = 1580.69642467072
$ws.Cells.Item(10, 19).Value = 0.006981952600720851
$ws.Cells.Item(10, 20).Value = 0.006981952600720851
$ws.Cells.Item(11, 7).Value = 51.27300266666666
$ws.Cells.Item(11, 8).Value = 153.819008
$ws.Cells.Item(11, 9).Value = 0.01048663484403512
$ws.Cells.Item(11, 10).Value = 0.01048663484403512
$ws.Cells.Item(11, 15).Value = 0.2094791321596951
$ws.Cells.Item(11, 16).Value = 0.2094791321596952
$ws.Cells.Item(11, 17).Value = 55.25937607799466
$ws.Cells.Item(11, 18).Value = 497.334384701952
$ws.Cells.Item(11, 19).Value = 0.002196731166404098
$ws.Cells.Item(11, 20).Value = 0.002196731166404098
$ws.Cells.Item(12, 7).Value = 51.27300266666666
$ws.Cells.Item(12, 8).Value = 153.819008
$ws.Cells.Item(12, 9).Value = 0.01048663484403512
$ws.Cells.Item(12, 10).Value = 0.01048663484403512
$ws.Cells.Item(12, 13).Value = 0.62317
$ws.Cells.Item(12, 14).Value = 1.86951
$ws.Cells.Item(12, 15).Value = 0.1211239647746572
$ws.Cells.Item(12, 16).Value = 0.1211239647746572
$ws.Cells.Item(12, 17).Value = 31.95179707178666
$ws.Cells.Item(12, 18).Value = 287.56617364608
$ws.Cells.Item(12, 19).Value = 0.001270182789453603
$ws.Cells.Item(12, 20).Value = 0.001270182789453603
$ws.Cells.Item(13, 7).Value = 51.27300266666666
$ws.Cells.Item(13, 8).Value = 153.819008
$ws.Cells.Item(13, 9).Value = 0.01048663484403512
$ws.Cells.Item(13, 10).Value = 0.01048663484403512
$ws.Cells.Item(13, 13).Value = 0.01852966666666667
$ws.Cells.Item(13, 14).Value = 0.055589
$ws.Cells.Item(13, 15).Value = 0.003601564087840353
$ws.Cells.Item(13, 16).Value = 0.003601564087840353
$ws.Cells.Item(13, 17).Value = 0.9500716484124444
$ws.Cells.Item(13, 18).Value = 8.550644835711999
$ws.Cells.Item(13, 19).Value = 0.00003776828745657223
$ws.Cells.Item(13, 20).Value = 0.00003776828745657223
$ws.Cells.Item(14, 7).Value = 2.837922333333333
$ws.Cells.Item(14, 8).Value = 8.513767
$ws.Cells.Item(14, 9).Value = 0.0005804273921477663
$ws.Cells.Item(14, 10).Value = 0.0005804273921477662
$ws.Cells.Item(14, 13).Value = 3.425446666666666
$ws.Cells.Item(14, 14).Value = 10.27634
$ws.Cells.Item(14, 15).Value = 0.6657953389778073
$ws.Cells.Item(14, 16).Value = 0.6657953389778073
$ws.Cells.Item(14, 17).Value = 9.721151596975554
$ws.Cells.Item(14, 18).Value = 87.49036437277999
$ws.Cells.Item(14, 19).Value = 0.0003864458523070267
$ws.Cells.Item(14, 20).Value = 0.0003864458523070267
$ws.Cells.Item(15, 7).Value = 2.837922333333333
$ws.Cells.Item(15, 8).Value = 8.513767
$ws.Cells.Item(15, 9).Value = 0.0005804273921477663
$ws.Cells.Item(15, 10).Value = 0.0005804273921477662
$ws.Cells.Item(15, 15).Value = 0.2094791321596951
$ws.Cells.Item(15, 16).Value = 0.2094791321596952
$ws.Cells.Item(15, 17).Value = 3.058565118905333
$ws.Cells.Item(15, 18).Value = 27.527086070148
$ws.Cells.Item(15, 19).Value = 0.0001215874263888291
$ws.Cells.Item(15, 20).Value = 0.0001215874263888291
$ws.Cells.Item(16, 7).Value = 2.837922333333333
$ws.Cells.Item(16, 8).Value = 8.513767
$ws.Cells.Item(16, 9).Value = 0.0005804273921477663
$ws.Cells.Item(16, 10).Value = 0.0005804273921477662
$ws.Cells.Item(16, 13).Value = 0.62317
$ws.Cells.Item(16, 14).Value = 1.86951
$ws.Cells.Item(16, 15).Value = 0.1211239647746572
$ws.Cells.Item(16, 16).Value = 0.1211239647746572
$ws.Cells.Item(16, 17).Value = 1.768508060463333
$ws.Cells.Item(16, 18).Value = 15.91657254417
$ws.Cells.Item(16, 19).Value = 0.00007030366700075218
$ws.Cells.Item(16, 20).Value = 0.00007030366700075216
$ws.Cells.Item(17, 7).Value = 2.837922333333333
$ws.Cells.Item(17, 8).Value = 8.513767
$ws.Cells.Item(17, 9).Value = 0.0005804273921477663
$ws.Cells.Item(17, 10).Value = 0.0005804273921477662
$ws.Cells.Item(17, 13).Value = 0.01852966666666667
$ws.Cells.Item(17, 14).Value = 0.055589
$ws.Cells.Item(17, 15).Value = 0.003601564087840353
$ws.Cells.Item(17, 16).Value = 0.003601564087840353
$ws.Cells.Item(17, 17).Value = 0.05258575486255555
$ws.Cells.Item(17, 18).Value = 0.473271793763
$ws.Cells.Item(17, 19).Value = 0.000002090446451158225
$ws.Cells.Item(17, 20).Value = 0.000002090446451158224
